$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the BAM-filter Cypher query text in cell B4 (shared string) ---
$newQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE  f.file_type in ['BAM']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

$ws.Cells.Item(4, 2).Value = $newQuery

# --- The longer query text now wraps across more lines, so the row grows taller ---
$ws.Rows.Item(4).RowHeight = 248

# --- Move the active selection from D2 to B5 ---
$ws.Range("B5").Select()
